# edit.ps1 - applies the "Starting again in 2018" edits:
#  1. Highlight the first "Stanford" mention and append trailing spaces
#     to the end of paragraph 1.
#  2. Remove the old _GoBack bookmark that sat after
#     "...random architecture."
#  3. Restructure "If admitted to Stanford..." paragraph: move the
#     _GoBack bookmark there (after "If "), and highlight "Stanford",
#     "Xiaoliang Qi's", second "Stanford", "Shoucheng Zhang's" and
#     "Steven Kivelson's".
#  4. Highlight "Stanford" in "...doctoral program at Stanford."
#  5. Split the header's "Stanford University" run so "Stanford" is
#     highlighted separately from " University".

$d = $word.ActiveDocument
$wdYellow = 7
$wdFindContinue = 1

function Get-SubRange {
    param(
        [__ComObject]$containerRange,
        [string]$needle
    )
    $sub = $d.Range($containerRange.Start, $containerRange.End)
    $found = $sub.Find.Execute($needle, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
    if (-not $found) {
        throw "Get-SubRange: '$needle' not found in container"
    }
    return $sub
}

function Highlight-InContext {
    param(
        [string]$context,
        [string]$needle
    )
    $ctx = $d.Content
    $found = $ctx.Find.Execute($context, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
    if (-not $found) {
        throw "Highlight-InContext: context '$context' not found"
    }
    $sub = Get-SubRange $ctx $needle
    $sub.Font.HighlightColorIndex = $wdYellow
    return $sub
}

# ---------------------------------------------------------------------
# 1. Paragraph 1: "...PhD program, and Stanford is an ideal location
#    for this study." -> highlight "Stanford" and add trailing spaces.
# ---------------------------------------------------------------------
Highlight-InContext "and Stanford is an ideal location for this study." "Stanford" | Out-Null

$endRng = $d.Content
$endRng.Find.Execute("for this study.", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$tailPos = $endRng.End
$tailRng = $d.Range($tailPos, $tailPos)
$tailRng.InsertAfter("   ")

# ---------------------------------------------------------------------
# 2. Remove the _GoBack bookmark after "...random architecture."
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 3. "If admitted to Stanford, ..." paragraph restructuring.
# ---------------------------------------------------------------------

# 3a. Insert the _GoBack bookmark right after "If ".
$ifRng = $d.Content
$ifRng.Find.Execute("If admitted to Stanford", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$bmPos = $ifRng.Start + 3
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3b. Highlight "Stanford" in "If admitted to Stanford,"
Highlight-InContext "If admitted to Stanford, I would be excited" "Stanford" | Out-Null

# 3c. Highlight "Xiaoliang" + " Qi's" (contiguous across the spell-check run split).
$xqCtx = $d.Content
$xqCtx.Find.Execute("Xiaoliang Qi’s research in quantum entanglement", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$xqSub = Get-SubRange $xqCtx "Xiaoliang Qi’s"
$xqSub.Font.HighlightColorIndex = $wdYellow

# 3d. Highlight second "Stanford" ("interested in at Stanford includes")
Highlight-InContext "interested in at Stanford includes" "Stanford" | Out-Null

# 3e. Highlight "Shoucheng" + " Zhang's" (contiguous).
$szCtx = $d.Content
$szCtx.Find.Execute("Shoucheng Zhang’s research", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$szSub = Get-SubRange $szCtx "Shoucheng Zhang’s"
$szSub.Font.HighlightColorIndex = $wdYellow

# 3f. Highlight "Steven " + "Kivelson's" (contiguous).
$skCtx = $d.Content
$skCtx.Find.Execute("Professor Steven Kivelson’s work", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$skSub = Get-SubRange $skCtx "Steven Kivelson’s"
$skSub.Font.HighlightColorIndex = $wdYellow

# ---------------------------------------------------------------------
# 4. "...doctoral program at Stanford." -> highlight "Stanford"
# ---------------------------------------------------------------------
Highlight-InContext "doctoral program at Stanford. I plan to pursue" "Stanford" | Out-Null

# ---------------------------------------------------------------------
# 5. Header: split "Stanford University" so "Stanford" is highlighted.
#    (Headers live in their own story range, so $d.Range(start,end)
#    -- which addresses the main-body story -- can't be used here;
#    narrow via a Duplicate of the header range instead.)
# ---------------------------------------------------------------------
$hdr = $d.Sections(1).Headers(1)
$hdrCtx = $hdr.Range
$hdrCtx.Find.Execute("Stanford University", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$hdrSub = $hdrCtx.Duplicate
$hdrSub.Find.Execute("Stanford", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$hdrSub.Font.HighlightColorIndex = $wdYellow

Write-Output "Done"
